# Fix SEA/SF name mix-up on Sheet1 and reset the active selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 29 is "San Francisco 49ers" - fff_abbreviation/vegas columns were
# incorrectly set to "SEA"; they should be "SF".
$ws.Range("C29").Value = "SF"
$ws.Range("D29").Value = "SF"

# Row 30 is "Seattle Seahawks" - fff_abbreviation/vegas columns were
# incorrectly set to "SF"; they should be "SEA".
$ws.Range("C30").Value = "SEA"
$ws.Range("D30").Value = "SEA"

# Reset the active cell/selection back to D1.
$ws.Range("D1").Select()
